# Daily attendance processing - 2025-12-28 12:49:02
# Normalize the "Recorded By" (column G) values: when the list of
# recorders starts with an exact "System" entry followed by at least one
# more entry, move "System" so it is no longer the first item (swap it
# with the entry that was originally second).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Value2

    if ($text -ne $null -and $text -ne "") {
        $parts = $text -split ", "

        if ($parts.Length -gt 1 -and $parts[0] -ceq "System") {
            $first = $parts[0]
            $second = $parts[1]
            $parts[0] = $second
            $parts[1] = $first
            $newText = [string]::Join(", ", $parts)

            if ($newText -ne $text) {
                $cell.Value = $newText
            }
        }
    }
}
